$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Replace all occurrences of the placeholder kid-name tokens with numeric IDs,
# across every cell in the sheet (handles both standalone "apple111"/"pear222"
# cells and "[apple111]"/"[pear222]" mentions embedded within sentences).
# LookAt:=2 (xlPart) so substrings inside larger sentences are matched too.
$ws.Cells.Replace("apple111", "id82", 2, 1, $false, $false, $false, $false)
$ws.Cells.Replace("pear222", "id83", 2, 1, $false, $false, $false, $false)

# Update the selected cell to match where the author ended up after editing.
$ws.Range("G145").Select()

$wb.Save()
